$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

$ws.Range("A2").Value = "iAU_TC_ID_159"
$ws.Range("B2").Value = "@RegressionA Validation of Edit an Exam ( Previleges: Applicable only for Exam Admin)"
